$wb = $excel.ActiveWorkbook
$poc = $wb.Worksheets.Item("poc_config")
$kpi = $wb.Worksheets.Item("kpi")

# --- Structural changes on the "kpi" sheet ---
# Insert 14 new columns (E..R) before the existing column E.
$kpi.Range("E1:R1").EntireColumn.Insert()
# Insert a new row before row 3 (keeps header row 1 / first data row 2 intact).
$kpi.Rows.Item(3).Insert()

# --- Copy formatting from "poc_config" so new cells pick up the exact
#     same styles already used elsewhere in the workbook ---
# Header look-alike cells G1:R1 reuse the poc_config sub-header style/text.
$poc.Range("F1:Q1").Copy()
$kpi.Range("G1:R1").PasteSpecial(-4163)
$poc.Range("F1:Q1").Copy()
$kpi.Range("G1:R1").PasteSpecial(-4122)

# Data rows 2:3 (A:R) reuse the bordered "data row" style from poc_config.
$poc.Range("A2").Copy()
$kpi.Range("A2:R3").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- New header text (row 1) ---
$kpi.Range("E1").Value = "Include Empty"
$kpi.Range("F1").Value = "Include Irrelevant"

# --- Row 2 values ---
$kpi.Range("A2").Value = "CCJP_RED_SCORE"
$kpi.Range("B2").Value = 0.16
$kpi.Range("C2").Value = 0.42
$kpi.Range("D2").Value = 0.42
$kpi.Range("E2").Value = "N/A"
$kpi.Range("F2").Value = "N/A"
$kpi.Range("G2").Value = "address_city"
$kpi.Range("H2").Value = "Tokyo"

# --- Row 3 values (new KPI row) ---
$kpi.Range("A3").Value = "CCJP_UNIQUE_DIST_OWN_MANU"
$kpi.Range("B3").Value = "N/A"
$kpi.Range("C3").Value = "N/A"
$kpi.Range("D3").Value = "N/A"
$kpi.Range("E3").Value = "N"
$kpi.Range("F3").Value = "N"
$kpi.Range("G3").Value = "address_city"
$kpi.Range("H3").Value = "Tokyo"

# --- Formatting tweaks shared by the header style (font 6 / style 4) ---
# Shrink the header font from 12pt to 8pt and give the header a border,
# matching the updated "kpi" banner look.
$kpi.Range("A1:F1").Font.Size = 8
$kpi.Range("A1:F1").Borders.LineStyle = 1

# --- Selection bookkeeping to mirror the saved cursor position ---
$poc.Range("D32").Select()
$kpi.Range("F11").Select()
